# Edit: "Tasks.xlsx" - mark Task #11 ("Removing expired products") as Finished,
# remove the now-redundant note "Same as 11 for expired products.", and update
# the active cell selection on the "Remaining Tasks" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Remaining Tasks")

# Row 15 (Task #11) - clone the formatting that is already used by other
# "Finished" rows (row 14) so the green status fill / borders match, then
# overwrite the cell values.
$ws.Range("A14:E14").Copy() | Out-Null
$ws.Range("A15:E15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(15, 2).Value = "Removing expired products"
$ws.Cells.Item(15, 3).Value = "Svetoslav"
$ws.Cells.Item(15, 4).Value = "Finished"
$ws.Cells.Item(15, 5).Value = ""

# Update the active cell / selection stored with the sheet view.
$ws.Range("I13").Select() | Out-Null

$wb.Save()
